$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A50").Value = 49
$ws.Range("B50").Value = 82
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 11
$ws.Range("E50").Value = 20
$ws.Range("F50").Value = 93
$ws.Range("G50").Value = 113
